$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 79, shifting existing rows 79..218 down to 80..219
$ws.Rows.Item(79).Insert()

# Populate the newly inserted row 79 with the new weekly data entry
$ws.Cells.Item(79,1).Value = 4
$ws.Cells.Item(79,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(79,3).Value = "Los Lagos"
$ws.Cells.Item(79,4).Value = 44533
$ws.Cells.Item(79,5).Value = 10
$ws.Cells.Item(79,6).Value = 100114014
$ws.Cells.Item(79,7).Value = "Betarraga"
$ws.Cells.Item(79,8).Value = "Sin especificar"
$ws.Cells.Item(79,9).Value = "Primera"
$ws.Cells.Item(79,10).Value = 1100
$ws.Cells.Item(79,11).Value = 1000
$ws.Cells.Item(79,12).Value = 1000
$ws.Cells.Item(79,13).Value = 1000
$ws.Cells.Item(79,14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(79,15).Value = "Región del Maule"
$ws.Cells.Item(79,16).Value = 200
$ws.Cells.Item(79,17).Value = 5
$ws.Cells.Item(79,18).Value = "Hortaliza"

# Match the date number format used by the rest of column D
$ws.Cells.Item(79,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
